$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells being updated hold plain text values (inlineStr), not real numbers,
# so force Text format first to keep Excel from reinterpreting the strings
# as numeric/percentage values (which would change their literal representation).
$textCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","E17","D18","E18","D19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D39","E39","D40","E40","D41","E42","D43","E43","E44","D45","E45","D46","E46","E47","D48","D50","E50","D51","E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "306.66"
$ws.Range("E2").Value = "-0.57%"
$ws.Range("D3").Value = "38.92"
$ws.Range("E3").Value = "7.06%"
$ws.Range("D4").Value = "5.096"
$ws.Range("E4").Value = "0.51%"
$ws.Range("D5").Value = "0.08072"
$ws.Range("E5").Value = "-0.50%"
$ws.Range("D6").Value = "1.956"
$ws.Range("E6").Value = "-4.06%"
$ws.Range("D7").Value = "4.191"
$ws.Range("E7").Value = "0.84%"
$ws.Range("D8").Value = "7.994"
$ws.Range("E8").Value = "1.70%"
$ws.Range("D9").Value = "0.9315"
$ws.Range("E9").Value = "0.36%"
$ws.Range("D10").Value = "0.1482"
$ws.Range("E10").Value = "3.52%"
$ws.Range("D11").Value = "0.1929"
$ws.Range("E11").Value = "-0.12%"
$ws.Range("D12").Value = "0.09141"
$ws.Range("E12").Value = "0.55%"
$ws.Range("D13").Value = "0.03496"
$ws.Range("E13").Value = "1.21%"
$ws.Range("D14").Value = "0.09782"
$ws.Range("E14").Value = "-1.32%"
$ws.Range("D15").Value = "0.001399"
$ws.Range("E15").Value = "-0.32%"
$ws.Range("D16").Value = "0.006060"
$ws.Range("E16").Value = "-2.89%"
$ws.Range("E17").Value = "-1.51%"
$ws.Range("D18").Value = "3.410"
$ws.Range("E18").Value = "-1.50%"
$ws.Range("D19").Value = "0.3424"
$ws.Range("D20").Value = "0.1303"
$ws.Range("E20").Value = "0.79%"
$ws.Range("D21").Value = "4.534"
$ws.Range("E21").Value = "-5.49%"
$ws.Range("D22").Value = "0.2414"
$ws.Range("E22").Value = "3.08%"
$ws.Range("D23").Value = "0.04370"
$ws.Range("E23").Value = "-0.43%"
$ws.Range("D24").Value = "0.001236"
$ws.Range("E24").Value = "0.25%"
$ws.Range("D25").Value = "0.004282"
$ws.Range("E25").Value = "-12.94%"
$ws.Range("D39").Value = "0.02037"
$ws.Range("E39").Value = "0.63%"
$ws.Range("D40").Value = "0.05101"
$ws.Range("E40").Value = "-1.14%"
$ws.Range("D41").Value = "0.007417"
$ws.Range("E42").Value = "1.41%"
$ws.Range("D43").Value = "0.1351"
$ws.Range("E43").Value = "-1.51%"
$ws.Range("E44").Value = "-0.90%"
$ws.Range("D45").Value = "0.009115"
$ws.Range("E45").Value = "-8.60%"
$ws.Range("D46").Value = "0.00006177"
$ws.Range("E46").Value = "-1.54%"
$ws.Range("E47").Value = "0.10%"
$ws.Range("D48").Value = "0.003101"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "0.10%"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "0.10%"
